# finger scintillator with silvering and rubber
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values ---
# ScintillatorThickness (row 2): 7 -> 10
$ws.Range("B2").Value = 10

# ScintillatorWidth (row 5): 145 -> 25
$ws.Range("B5").Value = 25

# ScintillatorHeight (row 6): 145 -> 200
$ws.Range("B6").Value = 200

# MountingScrewOffsett (row 13): 12.5 -> 9
$ws.Range("B13").Value = 9

# Recalculate dependent formulas (B7, B8) after the inputs above changed
$excel.CalculateFull()

# --- Add new rows for silvering & rubberized coating ---
$ws.Range("A15").Value = "silveringThickness"
$ws.Range("B15").Value = 0.2
$ws.Range("C15").Value = "mm"

$ws.Range("A16").Value = "rubberizedCoating"
$ws.Range("B16").Value = 0.5
$ws.Range("C16").Value = "mm"

# --- Update the active cell selection to A17 ---
$ws.Range("A17").Select()
